$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.708.10"
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").Value = "1.616.40"
$ws.Range("E3").Value = "  -1.62%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.36"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5070"
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2561"
$ws.Range("E8").Value = "  -1.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06345"
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.23"
$ws.Range("E10").Value = "  -3.35%  "
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.236"
$ws.Range("E12").Value = "  -1.90%  "
$ws.Range("D13").Value = "1.616.75"
$ws.Range("E13").Value = "  -1.87%  "
$ws.Range("D14").Value = "1.837.83"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5539"
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.39"
$ws.Range("E16").Value = "  -2.36%  "
$ws.Range("D17").Value = "0.0₅7488"
$ws.Range("E17").Value = "  -3.76%  "
$ws.Range("D18").Value = "25.711.34"
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.36"
$ws.Range("E20").Value = "  -3.31%  "
$ws.Range("E21").Value = "  -3.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.731"
$ws.Range("E22").Value = "  -3.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.952"
$ws.Range("E23").Value = "  -2.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.835"
$ws.Range("E25").Value = "  -3.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.47"
$ws.Range("E26").Value = "  -1.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1259"
$ws.Range("E27").Value = "  +3.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.712"
$ws.Range("E28").Value = "  -3.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.36"
$ws.Range("E29").Value = "  -2.45%  "
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04845"
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.285"
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.167"
$ws.Range("E33").Value = "  -2.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.544"
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.363"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8900"
$ws.Range("E36").Value = "  -3.32%  "
$ws.Range("D37").Value = "1.120.95"
$ws.Range("E37").Value = "  +0.59%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5477"
$ws.Range("E38").Value = "  -2.30%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.521"
$ws.Range("E39").Value = "  -3.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01554"
$ws.Range("E40").Value = "  -1.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9983"
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.562"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7900"
$ws.Range("E43").Value = "  -2.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.91"
$ws.Range("E44").Value = "  -3.02%  "
$ws.Range("D45").Value = "1.762.41"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("E46").Value = "  -5.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4402"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.51"
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05068"
$ws.Range("E49").Value = "  -3.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.490"
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9977"
$ws.Range("E51").Value = "  -0.66%  "
